# Update iServ stats for 2025-09 (row 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B22").Value = 6277
$ws.Range("C22").Value = 988
$ws.Range("D22").Value = 5703959
$ws.Range("E22").Value = 908.7078222080612
$ws.Range("F22").Value = 8.056464107419515
$ws.Range("G22").Value = 3.347280334728042
$ws.Range("H22").Value = 24.04333284040396
